$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ParticipantsTab row (row 2): replace the old, simplistic Participant-ID
# query (B2) with the corrected / expanded query that also threads through
# diagnosis / file / genomic_info and sorts + orders the results. ---
$newParticipantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Detection of Colorectal Cancer Susceptibility Loci Using Genome-Wide Sequencing"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@
# Normalize line endings to CRLF, matching the rest of the workbook, and
# drop the trailing newline the here-string literal adds.
$newParticipantsQuery = ($newParticipantsQuery -replace "`r?`n", "`r`n").TrimEnd("`r","`n")

$ws.Range("B2").Value2 = $newParticipantsQuery

# The cell now wraps across more lines, so the row needs to grow to fit.
$ws.Rows.Item(2).RowHeight = 299.25

# Selection moves from A2 to B2.
$ws.Range("B2").Select() | Out-Null
